$wb = $excel.ActiveWorkbook

# "展览" sheet (exhibitions) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1197
$ws1.Range("F4").Value = 14476
$ws1.Range("F5").Value = 17149
$ws1.Range("F9").Value = 51
$ws1.Range("F11").Value = 27
$ws1.Range("F12").Value = 53
$ws1.Range("F17").Value = 12
$ws1.Range("F18").Value = 125
$ws1.Range("F19").Value = 41
$ws1.Range("F20").Value = 1305
$ws1.Range("F21").Value = 145
$ws1.Range("F22").Value = 74
$ws1.Range("F23").Value = 59
$ws1.Range("F25").Value = 7063
$ws1.Range("F27").Value = 31
$ws1.Range("F28").Value = 1154
$ws1.Range("F29").Value = 28
$ws1.Range("F31").Value = 43
$ws1.Range("F32").Value = 5823
$ws1.Range("F33").Value = 129
$ws1.Range("F35").Value = 216
$ws1.Range("F36").Value = 4989
$ws1.Range("F37").Value = 29

# "全部类型" sheet (all types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1197
$ws4.Range("F4").Value = 14476
$ws4.Range("F5").Value = 17149
$ws4.Range("F9").Value = 51
$ws4.Range("F11").Value = 27
$ws4.Range("F12").Value = 53
$ws4.Range("F17").Value = 12
$ws4.Range("F18").Value = 125
$ws4.Range("F19").Value = 41
$ws4.Range("F20").Value = 1305
$ws4.Range("F21").Value = 145
$ws4.Range("F22").Value = 74
$ws4.Range("F24").Value = 59
$ws4.Range("F26").Value = 7063
$ws4.Range("F28").Value = 31
$ws4.Range("F29").Value = 1154
$ws4.Range("F30").Value = 28
$ws4.Range("F32").Value = 43
$ws4.Range("F34").Value = 5823
$ws4.Range("F35").Value = 129
$ws4.Range("F37").Value = 216
$ws4.Range("F38").Value = 4989
$ws4.Range("F39").Value = 29
